$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the old trailing rows (8-14); shifts nothing else, just drops them ---
$ws.Range("A8:A14").EntireRow.Delete()

# --- Row 1: rename existing header, add three new headers ---
$ws.Range("B1").Value = "CENTRAL"
$ws.Range("C1").Value = "X1"
$ws.Range("D1").Value = "X2"
$ws.Range("E1").Value = "X3"

# Give the new header cells the same look (font/border/alignment) as B1,
# without touching the shared style table beyond what's already there.
$ws.Range("B1").Copy()
$ws.Range("C1:E1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Rows 2-6: drop the old single "B" value column, populate C/D/E ---
$ws.Range("B2:B6").ClearContents()

$ws.Range("C2").Value = 14.10361184023227
$ws.Range("D2").Value = 31.88389669424919
$ws.Range("E2").Value = 74.65580668962787

$ws.Range("C3").Value = 25.14925804642162
$ws.Range("D3").Value = 37.0349530383237
$ws.Range("E3").Value = 55.45427847716447

$ws.Range("C4").Value = 42.39742708169723
$ws.Range("D4").Value = 42.39742708169723
$ws.Range("E4").Value = 42.39742708169723

$ws.Range("C5").Value = 71.47494450253166
$ws.Range("D5").Value = 48.53636026722488
$ws.Range("E5").Value = 32.41484467042599

$ws.Range("C6").Value = 127.452587571939
$ws.Range("D6").Value = 56.37773326094273
$ws.Range("E6").Value = 24.07772285712328

# --- Row 7: turn the old numeric index (6) into the "CENTRAL" label,
#     and set B7 to the diagonal value from row 4 ---
$ws.Range("A7").Value = "CENTRAL"
$ws.Range("B7").Value = 42.39742708169723
